# Apply the authored changes to "WhatsUrSay_Deliverable 1.pptx":
#  1. Slide 3 ("What & Why?"): fix "specific group or users." -> "specific group of users."
#  2. Slide 3: merge the split runs of the "WhatUrSay would provide..." sentence
#     back into a single run (no textual change - just a run-merge).
#  3. Slide 8 ("Scope of the Project (cont..)"): grow the table's cached overall
#     height (p:xfrm ext cy) from 4624121 EMU to 4661458 EMU.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1 & 2. Slide 3 - "Content Placeholder 2" text box
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$contentShape = $slide3.Shapes.Item("Content Placeholder 2")
$tr = $contentShape.TextFrame.TextRange

# --- Edit 1: "...specific group or users." -> "...specific group of users." ---
$pollsPara = $tr.Paragraphs(3)
$offset = $pollsPara.Text.IndexOf("group or users")
if ($offset -ge 0) {
    $wordRange = $tr.Characters($pollsPara.Start + $offset, 9)
    $wordRange.Text = "group of "
}

# --- Edit 2: merge the 3 runs making up the "WhatUrSay would provide..." ------
#     sentence into a single run (text itself is unchanged).
$sayPara = $tr.Paragraphs(8)
$prefix = "WhatUrSay"
$tailRange = $tr.Characters($sayPara.Start + $prefix.Length, $sayPara.Length - $prefix.Length)
$tailRange.Text = " would provide an easy and user friendly way to create, participate and result publishing of polls/surveys."

# ---------------------------------------------------------------------------
# 3. Slide 8 - grow the table's cached overall extent
# ---------------------------------------------------------------------------
# The graphicFrame's cached <a:ext cy="..."/> is recomputed (as the sum of all
# row heights) as soon as any row's Height is touched, and the rendered Height
# getter does not always echo back the exact stored EMU value for a row - so
# the new height for the (touched) last row is derived from the other,
# untouched rows' heights rather than from a relative delta off the getter.
$slide8 = $p.Slides.Item(8)
$tableShape = $slide8.Shapes.Item("Content Placeholder 3")
$table = $tableShape.Table

$emuPerPoint = 12700
$otherRowsHeightEmu = 0
for ($i = 1; $i -lt $table.Rows.Count; $i++) {
    $otherRowsHeightEmu += [Math]::Round($table.Rows.Item($i).Height * $emuPerPoint)
}

$targetHeightEmu = 4661458
$lastRow = $table.Rows.Item($table.Rows.Count)
$lastRow.Height = ($targetHeightEmu - $otherRowsHeightEmu) / $emuPerPoint
